$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 2.6
$ws.Range("G2").Value = 2.64
$ws.Range("H2").Value = 2.9
$ws.Range("I2").Value = 2.94
$ws.Range("J2").Value = 3.65
$ws.Range("K2").Value = 3.7
$ws.Range("L2").Value = 1.33
$ws.Range("N2").Value = 3.75
$ws.Range("O2").Value = 1.32
$ws.Range("P2").Value = 1.96
$ws.Range("Q2").Value = 1.94
$ws.Range("S2").Value = 3.4
$ws.Range("T2").Value = 1.74
$ws.Range("W2").Value = 1.61
$ws.Range("Y2").Value = 14.5
$ws.Range("AF2").Value = 20
$ws.Range("AI2").Value = 48
$ws.Range("AL2").Value = 46
$ws.Range("H3").Value = 2.2
$ws.Range("I3").Value = 2.24
$ws.Range("N3").Value = 7
$ws.Range("P3").Value = 2.98
$ws.Range("R3").Value = 1.83
$ws.Range("S3").Value = 2.14
$ws.Range("X3").Value = 38
$ws.Range("Y3").Value = 20
$ws.Range("Z3").Value = 22
$ws.Range("AB3").Value = 27
$ws.Range("AF3").Value = 36
$ws.Range("AG3").Value = 18
$ws.Range("AH3").Value = 15
$ws.Range("AI3").Value = 30
$ws.Range("AK3").Value = 36
$ws.Range("AL3").Value = 38
$ws.Range("AM3").Value = 55
$ws.Range("AN3").Value = 19
$ws.Range("J5").Value = 3.25
$ws.Range("P5").Value = 1.76
$ws.Range("J6").Value = 3.85
$ws.Range("R7").Value = 1.35
$ws.Range("U7").Value = 1.93
$ws.Range("AJ7").Value = 18.5
$ws.Range("T8").Value = 2.32
$ws.Range("G9").Value = 3.95
$ws.Range("J9").Value = 3.35
$ws.Range("P10").Value = 1.87
$ws.Range("Q10").Value = 2
